$d = $word.ActiveDocument

$texts = @(
  "במחצית זאת למדנו חומש שמות, הרחבנו על עשרת המכות בצורה חווייתית,הבאנו מדרשים ומשלים על הנ`"ל.`nעדי את תלמידה מדהימה, בהצלחה!",
  "במחצית זאת הכרנו את אותיות האנגלית, למדנו את אותיות הניקוד, והרחבנו על כל אות בנפרד, עשינו זאת כחוויה, והייתה אוירה טובה,`nעדי את תלמידה מקסימה, בהצלחה!",
  "במחצית זאת למדנו את תורת המספרים,הכרנו את הפעולות הבסיסיות, חיבור וחיסור, התעסקנו עם מספרים גבוהים יותר,וניסנו לעשות את הפעולות הבסיסיות עליהם.`nעדי את ילדה נהדרת, בהצלחה רבה!",
  "במחצית זאת למדנו על מעגל השנה,על כל חג הרחבנו ופירטנו את מקומינו ביחס לחג, הייתה אוירה מרוממת.`nעדי אתה בחור עם שאיפות, הרבה הצלחה!"
)

for ($i = 1; $i -le $d.Tables.Count; $i++) {
  $tbl = $d.Tables.Item($i)
  $cell = $tbl.Cell(1, 2)
  $cell.Range.Text = $texts[$i - 1]
}

